$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = -7.789999999999999
$ws.Range("D7").Value = -7.114999999999999
$ws.Range("C8").Value = -12.672
$ws.Range("B12").Value = 5.513
$ws.Range("C12").Value = -13.073
$ws.Range("C14").Value = -11.675
$ws.Range("D19").Value = -7.663999999999999
$ws.Range("D21").Value = -7.597999999999999
$ws.Range("C22").Value = -12.846
$ws.Range("D24").Value = -7.637
